$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 50.5
$ws.Range("I5").Value = 50.5
$ws.Range("K5").Value = 50.5
$ws.Range("M5").Value = 64.5

$ws.Range("H28").Value = 528.75
$ws.Range("I28").Value = 528.75
$ws.Range("K28").Value = 528.75
$ws.Range("M28").Value = -43.75

$ws.Range("H40").Value = 6859.375
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 8360
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 8360
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -8710

$ws.Range("H86").Value = 3309.6
$ws.Range("I86").Value = 1900
$ws.Range("K86").Value = 1900
$ws.Range("M86").Value = -777

$ws.Range("H89").Value = 3309.6
$ws.Range("I89").Value = 1900
$ws.Range("K89").Value = 9500
$ws.Range("M89").Value = -3884

$ws.Range("H132").Value = 3322.7144
$ws.Range("I132").Value = 651.8
$ws.Range("K132").Value = 1955.4
$ws.Range("M132").Value = 574.6000000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1933
$ws.Range("I45").Value = 1933
$ws.Range("K45").Value = 1933
$ws.Range("M45").Value = -1556

$ws.Range("H88").Value = 368.8
$ws.Range("I88").Value = 6
$ws.Range("J88").Value = 459.5
$ws.Range("K88").Value = 6
$ws.Range("L88").Value = 459.5
$ws.Range("M88").Value = 400
$ws.Range("N88").Value = -1271.5

$ws.Range("H91").Value = 368.8
$ws.Range("I91").Value = 6
$ws.Range("J91").Value = 459.5
$ws.Range("K91").Value = 6
$ws.Range("L91").Value = 459.5
$ws.Range("M91").Value = 1398
$ws.Range("N91").Value = -3267.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1986
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 1986
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 331.25
$ws.Range("I7").Value = 187.5
$ws.Range("K7").Value = 187.5
$ws.Range("M7").Value = -74.5

$ws.Range("H94").Value = 1724
$ws.Range("I94").Value = 1698.5
$ws.Range("K94").Value = 1698.5
$ws.Range("M94").Value = -1247.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 801.6
$ws.Range("J34").Value = 967.6667
$ws.Range("L34").Value = 2903.0001
$ws.Range("N34").Value = -3071.0001

$ws.Range("H39").Value = 4700
$ws.Range("I39").Value = 4500
$ws.Range("K39").Value = 13500
$ws.Range("M39").Value = -13206

$ws.Range("H55").Value = 1951.5454
$ws.Range("I55").Value = 1266.7273
$ws.Range("J55").Value = 2636.3635
$ws.Range("K55").Value = 3800.1819
$ws.Range("L55").Value = 7909.0905
$ws.Range("M55").Value = -3623.1819
$ws.Range("N55").Value = -8263.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 20000

$ws.Range("H132").Value = 2800.6
$ws.Range("I132").Value = 2071.5
$ws.Range("J132").Value = 4258.8
$ws.Range("K132").Value = 6214.5
$ws.Range("L132").Value = 12776.4
$ws.Range("M132").Value = -3684.5
$ws.Range("N132").Value = -17836.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5250
$ws.Range("I7").Value = 5250
$ws.Range("K7").Value = 5250
$ws.Range("M7").Value = -5138

$ws.Range("H16").Value = 1999.5
$ws.Range("I16").Value = 1999.5
$ws.Range("K16").Value = 1999.5
$ws.Range("M16").Value = -1829.5

$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H40").Value = 5988.846
$ws.Range("I40").Value = 5904.5835
$ws.Range("K40").Value = 5904.5835
$ws.Range("M40").Value = -5768.5835

$ws.Range("H46").Value = 3250.6
$ws.Range("I46").Value = 1332.6666
$ws.Range("K46").Value = 1332.6666
$ws.Range("M46").Value = -1144.6666

$ws.Range("H68").Value = 2467.8333
$ws.Range("I68").Value = 2451
$ws.Range("J68").Value = 2501.5
$ws.Range("K68").Value = 2451
$ws.Range("L68").Value = 2501.5
$ws.Range("M68").Value = -1702
$ws.Range("N68").Value = -3999.5

$ws.Range("H71").Value = 2467.8333
$ws.Range("I71").Value = 2451
$ws.Range("J71").Value = 2501.5
$ws.Range("K71").Value = 12255
$ws.Range("L71").Value = 12507.5
$ws.Range("M71").Value = -8511
$ws.Range("N71").Value = -19995.5

$ws.Range("H82").Value = 1841.6666
$ws.Range("I82").Value = 1730
$ws.Range("K82").Value = 1730
$ws.Range("M82").Value = -1369

$ws.Range("H85").Value = 1841.6666
$ws.Range("I85").Value = 1730
$ws.Range("K85").Value = 1730
$ws.Range("M85").Value = -482

$ws.Range("H122").Value = 4899.8
$ws.Range("J122").Value = 5500
$ws.Range("L122").Value = 16500
$ws.Range("N122").Value = -21400

$ws.Range("H126").Value = 5250
$ws.Range("I126").Value = 5250
$ws.Range("K126").Value = 15750
$ws.Range("M126").Value = -13280

$ws.Range("H130").Value = 10000
$ws.Range("J130").Value = 10000
$ws.Range("L130").Value = 10000
$ws.Range("N130").Value = -20040

$ws.Range("H132").Value = 8450
$ws.Range("I132").Value = 8450
$ws.Range("K132").Value = 25350
$ws.Range("M132").Value = -22820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1988.3
$ws.Range("I107").Value = 1431.4445
$ws.Range("K107").Value = 4294.333500000001
$ws.Range("M107").Value = -2374.333500000001

$ws.Range("H113").Value = 604.375
$ws.Range("I113").Value = 606.3333
$ws.Range("K113").Value = 1818.9999
$ws.Range("M113").Value = 351.0001

$ws.Range("H132").Value = 4490.2
$ws.Range("I132").Value = 2356.8572
$ws.Range("K132").Value = 7070.571599999999
$ws.Range("M132").Value = -4540.571599999999

$ws.Range("H136").Value = 875
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 750
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 2250
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -7350

